$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 412-470: columns D, J, K, L, M, N, P change; columns A,B,C,E,F,G,H,I,O,Q,R stay the same.
$data = @(
    ,@(412, 45124, 60, 19000, 19000, 19000, '$/caja 10 kilos', 1900)
    ,@(413, 45124, 50, 20000, 20000, 20000, '$/malla 10 kilos', 2000)
    ,@(414, 44902, 50, 17000, 17000, 17000, '$/caja 10 kilos', 1700)
    ,@(415, 44902, 50, 18000, 18000, 18000, '$/malla 10 kilos', 1800)
    ,@(416, 44631, 60, 19000, 20000, 19500, '$/caja 10 kilos', 1950)
    ,@(417, 44554, 100, 20000, 21000, 20500, '$/malla 10 kilos', 2050)
    ,@(418, 45070, 50, 18000, 18000, 18000, '$/caja 10 kilos', 1800)
    ,@(419, 45070, 50, 19000, 19000, 19000, '$/malla 10 kilos', 1900)
    ,@(420, 44790, 80, 27000, 28000, 27500, '$/malla 10 kilos', 2750)
    ,@(421, 44260, 70, 12000, 13000, 12571, '$/caja 10 kilos', 1257)
    ,@(422, 44272, 70, 12000, 12500, 12357, '$/caja 10 kilos', 1236)
    ,@(423, 45040, 60, 17000, 18000, 17500, '$/caja 10 kilos', 1750)
    ,@(424, 45040, 60, 19000, 20000, 19500, '$/malla 10 kilos', 1950)
    ,@(425, 44826, 60, 23000, 24000, 23500, '$/caja 10 kilos', 2350)
    ,@(426, 44692, 100, 20000, 21000, 20500, '$/caja 10 kilos', 2050)
    ,@(427, 44589, 100, 18000, 19000, 18500, '$/caja 10 kilos', 1850)
    ,@(428, 44565, 80, 19000, 20000, 19500, '$/caja 10 kilos', 1950)
    ,@(429, 44901, 50, 18000, 18000, 18000, '$/malla 10 kilos', 1800)
    ,@(430, 44322, 120, 14000, 15000, 14500, '$/caja 10 kilos', 1450)
    ,@(431, 44495, 100, 16000, 17000, 16500, '$/caja 10 kilos', 1650)
    ,@(432, 45111, 60, 20000, 21000, 20500, '$/malla 10 kilos', 2050)
    ,@(433, 44417, 120, 14000, 15000, 14500, '$/caja 10 kilos', 1450)
    ,@(434, 44323, 40, 15000, 16000, 15500, '$/caja 10 kilos', 1550)
    ,@(435, 44221, 65, 12500, 13000, 12731, '$/caja 10 kilos', 1273)
    ,@(436, 44664, 100, 20000, 21000, 20500, '$/caja 10 kilos', 2050)
    ,@(437, 44511, 100, 15000, 16000, 15500, '$/caja 10 kilos', 1550)
    ,@(438, 44253, 70, 13000, 13500, 13214, '$/caja 10 kilos', 1321)
    ,@(439, 44358, 120, 13000, 14000, 13500, '$/caja 10 kilos', 1350)
    ,@(440, 44235, 80, 12000, 13000, 12500, '$/caja 10 kilos', 1250)
    ,@(441, 44634, 60, 19000, 20000, 19500, '$/caja 10 kilos', 1950)
    ,@(442, 44420, 120, 14000, 15000, 14500, '$/caja 10 kilos', 1450)
    ,@(443, 44924, 60, 17000, 18000, 17500, '$/malla 10 kilos', 1750)
    ,@(444, 44638, 60, 18000, 19000, 18500, '$/caja 10 kilos', 1850)
    ,@(445, 44239, 40, 12500, 13000, 12750, '$/caja 10 kilos', 1275)
    ,@(446, 45075, 50, 17000, 17000, 17000, '$/caja 10 kilos', 1700)
    ,@(447, 45075, 50, 18000, 18000, 18000, '$/malla 10 kilos', 1800)
    ,@(448, 44971, 60, 17000, 18000, 17500, '$/caja 10 kilos', 1750)
    ,@(449, 44364, 120, 14000, 15000, 14500, '$/caja 10 kilos', 1450)
    ,@(450, 44517, 100, 18000, 19000, 18500, '$/caja 10 kilos', 1850)
    ,@(451, 44985, 50, 21000, 21000, 21000, '$/malla 10 kilos', 2100)
    ,@(452, 44859, 60, 16000, 17000, 16500, '$/malla 10 kilos', 1650)
    ,@(453, 44811, 60, 23000, 24000, 23500, '$/caja 10 kilos', 2350)
    ,@(454, 44637, 60, 19000, 20000, 19500, '$/caja 10 kilos', 1950)
    ,@(455, 45112, 60, 18000, 19000, 18500, '$/caja 10 kilos', 1850)
    ,@(456, 45112, 50, 20000, 20000, 20000, '$/malla 10 kilos', 2000)
    ,@(457, 45112, 60, 20000, 21000, 20500, '$/malla 10 kilos', 2050)
    ,@(458, 44802, 30, 27000, 28000, 27500, '$/caja 10 kilos', 2750)
    ,@(459, 44574, 100, 20000, 21000, 20500, '$/caja 10 kilos', 2050)
    ,@(460, 44594, 60, 18000, 19000, 18500, '$/caja 10 kilos', 1850)
    ,@(461, 44749, 80, 20000, 21000, 20500, '$/caja 10 kilos', 2050)
    ,@(462, 44957, 60, 18000, 19000, 18500, '$/malla 10 kilos', 1850)
    ,@(463, 45089, 50, 18000, 18000, 18000, '$/caja 10 kilos', 1800)
    ,@(464, 45089, 60, 19000, 19000, 19000, '$/malla 10 kilos', 1900)
    ,@(465, 45099, 40, 18000, 18000, 18000, '$/caja 10 kilos', 1800)
    ,@(466, 45099, 40, 21000, 21000, 21000, '$/malla 10 kilos', 2100)
    ,@(467, 45121, 40, 19000, 19000, 19000, '$/caja 10 kilos', 1900)
    ,@(468, 45121, 30, 21000, 21000, 21000, '$/malla 10 kilos', 2100)
    ,@(469, 44341, 60, 14000, 15000, 14500, '$/caja 10 kilos', 1450)
    ,@(470, 44662, 100, 19000, 20000, 19500, '$/caja 10 kilos', 1950)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Range("D$r").Value = $entry[1]
    $ws.Range("J$r").Value = $entry[2]
    $ws.Range("K$r").Value = $entry[3]
    $ws.Range("L$r").Value = $entry[4]
    $ws.Range("M$r").Value = $entry[5]
    $ws.Range("N$r").Value = $entry[6]
    $ws.Range("P$r").Value = $entry[7]
}

# Add new rows 471 and 472 (full rows, shifted from old rows 469 and 470).
# Row 471
$ws.Range("A471").Value = 7
$ws.Range("D471").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B471").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C471").Value = 'Ñuble'
$ws.Range("D471").Value = 44607
$ws.Range("E471").Value = 16
$ws.Range("F471").Value = 100112003
$ws.Range("G471").Value = 'Ajo'
$ws.Range("H471").Value = 'Chino'
$ws.Range("I471").Value = 'Primera'
$ws.Range("J471").Value = 60
$ws.Range("K471").Value = 19000
$ws.Range("L471").Value = 20000
$ws.Range("M471").Value = 19500
$ws.Range("N471").Value = '$/caja 10 kilos'
$ws.Range("O471").Value = 'China'
$ws.Range("P471").Value = 1950
$ws.Range("Q471").Value = 10
$ws.Range("R471").Value = 'Hortaliza'

# Row 472
$ws.Range("A472").Value = 7
$ws.Range("D472").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B472").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C472").Value = 'Ñuble'
$ws.Range("D472").Value = 45072
$ws.Range("E472").Value = 16
$ws.Range("F472").Value = 100112003
$ws.Range("G472").Value = 'Ajo'
$ws.Range("H472").Value = 'Chino'
$ws.Range("I472").Value = 'Primera'
$ws.Range("J472").Value = 60
$ws.Range("K472").Value = 17000
$ws.Range("L472").Value = 18000
$ws.Range("M472").Value = 17500
$ws.Range("N472").Value = '$/caja 10 kilos'
$ws.Range("O472").Value = 'China'
$ws.Range("P472").Value = 1750
$ws.Range("Q472").Value = 10
$ws.Range("R472").Value = 'Hortaliza'
